# Apply the "complete 15-shipment dataset" update to the Sample_Data workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a string value into a cell while guaranteeing it is stored as
# text (prevents Excel's automatic number/date inference for things like
# postal codes "00100" or ISO dates "2024-11-11"), and afterwards strip the
# temporary number-format override so the cell is left with no explicit
# style, matching plain data cells elsewhere in the sheet.
# ---------------------------------------------------------------------------
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# ---------------------------------------------------------------------------
# Rename the worksheet/tab
# ---------------------------------------------------------------------------
$ws.Name = "Shipments"

# ---------------------------------------------------------------------------
# Header row updates (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Unique identifier"
$ws.Range("C1").Value = "Collection Area Name"
$ws.Range("D1").Value = "Delivery Area Name"
$ws.Range("H1").Value = "Shipper"
$ws.Range("I1").Value = "Emission type by shipper mode"
$ws.Range("J1").Value = "Emission factor of emission type"

# K1/L1 are brand-new header cells; give them the same bold/centered header
# style as the rest of row 1 (copy format from A1) before setting their text.
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Weight"
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("L1").Value = "Cost of shipment"

# ---------------------------------------------------------------------------
# Existing rows 2-11: the old "Weight (kg)"/"Cost" columns (I, J) shift right
# to become columns K and L, and two new columns are inserted at I (emission
# type by shipper mode) and J (emission factor of emission type).
# ---------------------------------------------------------------------------
$existingShipperModes = @{
    2  = @{ Mode = "Van";             Factor = 0.21 }
    3  = @{ Mode = "Truck";           Factor = 0.18 }
    4  = @{ Mode = "Air Freight";     Factor = 0.65 }
    5  = @{ Mode = "Truck";           Factor = 0.22 }
    6  = @{ Mode = "Air Freight";     Factor = 0.62 }
    7  = @{ Mode = "Van";             Factor = 0.19 }
    8  = @{ Mode = "Container Ship";  Factor = 0.01 }
    9  = @{ Mode = "Truck";           Factor = 0.23 }
    10 = @{ Mode = "Air Freight";     Factor = 0.68 }
    11 = @{ Mode = "Van";             Factor = 0.20 }
}

foreach ($r in 2..11) {
    $info = $existingShipperModes[$r]

    # Move old Weight (I) and Cost (J) values into K and L first, before
    # overwriting column I with the new "shipper mode" text.
    $oldWeight = $ws.Range("I$r").Value()
    $oldCost   = $ws.Range("J$r").Value()

    $ws.Range("K$r").Value = $oldWeight
    $ws.Range("L$r").Value = $oldCost

    $ws.Range("I$r").Value = $info.Mode
    $ws.Range("J$r").Value = $info.Factor
}

# ---------------------------------------------------------------------------
# New rows 12-16 (SHIP011 - SHIP015)
# ---------------------------------------------------------------------------

# Row 12 - SHIP011
Set-TextValue "A12" "SHIP011"
Set-TextValue "B12" "2024-11-11"
Set-TextValue "C12" "2000"
Set-TextValue "D12" "3000"
$ws.Range("E12").Value = "Australia"
$ws.Range("F12").Value = "Australia"
$ws.Range("G12").Value = "Road"
$ws.Range("H12").Value = "Australia Post"
$ws.Range("I12").Value = "Truck"
$ws.Range("J12").Value = 0.25
$ws.Range("K12").Value = 20
$ws.Range("L12").Value = 35

# Row 13 - SHIP012
Set-TextValue "A13" "SHIP012"
Set-TextValue "B13" "2024-11-12"
$ws.Range("C13").Value = "W1A 1AA"
$ws.Range("D13").Value = "100-0001"
$ws.Range("E13").Value = "United Kingdom"
$ws.Range("F13").Value = "Japan"
$ws.Range("G13").Value = "Air"
$ws.Range("H13").Value = "FedEx"
$ws.Range("I13").Value = "Air Freight"
$ws.Range("J13").Value = 0.7
$ws.Range("K13").Value = 7.5
$ws.Range("L13").Value = 145

# Row 14 - SHIP013
Set-TextValue "A14" "SHIP013"
Set-TextValue "B14" "2024-11-13"
$ws.Range("C14").Value = "OX1 2JD"
$ws.Range("D14").Value = "CB2 1TN"
$ws.Range("E14").Value = "United Kingdom"
$ws.Range("F14").Value = "United Kingdom"
$ws.Range("G14").Value = "Road"
$ws.Range("H14").Value = "DPD"
$ws.Range("I14").Value = "Van"
$ws.Range("J14").Value = 0.21
$ws.Range("K14").Value = 4
$ws.Range("L14").Value = 9

# Row 15 - SHIP014
Set-TextValue "A15" "SHIP014"
Set-TextValue "B15" "2024-11-14"
Set-TextValue "C15" "75001"
Set-TextValue "D15" "00100"
$ws.Range("E15").Value = "France"
$ws.Range("F15").Value = "Italy"
$ws.Range("G15").Value = "Road"
$ws.Range("H15").Value = "TNT"
$ws.Range("I15").Value = "Truck"
$ws.Range("J15").Value = 0.24
$ws.Range("K15").Value = 30
$ws.Range("L15").Value = 65

# Row 16 - SHIP015
Set-TextValue "A16" "SHIP015"
Set-TextValue "B16" "2024-11-15"
$ws.Range("C16").Value = "CR0 1EA"
$ws.Range("D16").Value = "BN1 1AL"
$ws.Range("E16").Value = "United Kingdom"
$ws.Range("F16").Value = "United Kingdom"
$ws.Range("G16").Value = "Road"
$ws.Range("H16").Value = "Royal Mail"
$ws.Range("I16").Value = "Van"
$ws.Range("J16").Value = 0.19
$ws.Range("K16").Value = 2.5
$ws.Range("L16").Value = 6.5
